$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 16: "category" style A16 (reuse existing text "# Python", shared string index 25)
# plus the two brand-new entries for the "Range usage" topic.
$ws.Range("A16").Value = "# Python"
$ws.Range("B16").Value = "## Range usage"
$ws.Range("C16").Value = ">>> for i in range(5):    //print 0,1,2,3,4`n>>> for i in range(3, 6):    //print 3,4,5`n>>> for i in range(4, 10, 2):    //print 4,6,8`n>>> for i in range(0, -10, -2):    //print 0,-2,-4,-6,-8"

# Match the formatting already used by the rest of the "Category" (col A) /
# "Topic" + "Details" (cols B/C, word-wrapped) columns for this table.
$ws.Range("C16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 51

# Reflect the new selection state left behind by the edit.
$ws.Range("C22:C29").Select()
